# The presentation's theme (ppt/theme/theme1.xml, used by the slide
# master / all slides & layouts) is switched from the custom "Integral"
# colour palette to the standard Office colour palette:
#   dk2  : 455F51 -> 44546A
#   lt2  : E3DED1 -> E7E6E6
#   acc1 : 99CB38 -> 5B9BD5
#   acc2 : 63A537 -> ED7D31
#   acc3 : E6D024 -> A5A5A5
#   acc4 : CC9700 -> FFC000
#   acc5 : 4EB3CF -> 4472C4
#   acc6 : 378DA6 -> 70AD47
#   hlink: 6B9F25 -> 0563C1
#   folHl: B26B02 -> 954F72
# (dk1/lt1 stay black/white in both palettes.)
#
# PowerPoint exposes the theme's colour scheme slots through
# ThemeColorScheme.Item(1..12), in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# The .RGB setter takes a COM colour long (0x00BBGGRR).

function HexToCOMColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# The presentation has a single Design; its SlideMaster's Theme backs
# ppt/theme/theme1.xml.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = HexToCOMColor "000000"   # dk1
$colorScheme.Item(2).RGB  = HexToCOMColor "FFFFFF"   # lt1
$colorScheme.Item(3).RGB  = HexToCOMColor "44546A"   # dk2
$colorScheme.Item(4).RGB  = HexToCOMColor "E7E6E6"   # lt2
$colorScheme.Item(5).RGB  = HexToCOMColor "5B9BD5"   # accent1
$colorScheme.Item(6).RGB  = HexToCOMColor "ED7D31"   # accent2
$colorScheme.Item(7).RGB  = HexToCOMColor "A5A5A5"   # accent3
$colorScheme.Item(8).RGB  = HexToCOMColor "FFC000"   # accent4
$colorScheme.Item(9).RGB  = HexToCOMColor "4472C4"   # accent5
$colorScheme.Item(10).RGB = HexToCOMColor "70AD47"   # accent6
$colorScheme.Item(11).RGB = HexToCOMColor "0563C1"   # hlink
$colorScheme.Item(12).RGB = HexToCOMColor "954F72"   # folHlink
